$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '28.334.46'
$ws.Cells.Item(2, 5).Value = '  +4.05%  '
$ws.Cells.Item(3, 4).Value = '1.730.13'
$ws.Cells.Item(3, 5).Value = '  +2.38%  '
$ws.Cells.Item(4, 5).Value = '  -0.14%  '
$ws.Cells.Item(5, 4).Value = '''219.35'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.38%  '
$ws.Cells.Item(6, 5).Value = '  +0.39%  '
$ws.Cells.Item(7, 5).Value = '  -0.17%  '
$ws.Cells.Item(8, 4).Value = '''24.07'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +3.68%  '
$ws.Cells.Item(9, 5).Value = '  +2.05%  '
$ws.Cells.Item(10, 4).Value = '''0.0637'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +1.37%  '
$ws.Cells.Item(11, 4).Value = '''0.0894'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.26%  '
$ws.Cells.Item(12, 4).Value = '1.975.59'
$ws.Cells.Item(12, 5).Value = '  +2.45%  '
$ws.Cells.Item(13, 4).Value = '1.728.39'
$ws.Cells.Item(13, 5).Value = '  +2.04%  '
$ws.Cells.Item(14, 5).Value = '  +1.05%  '
$ws.Cells.Item(15, 4).Value = '''0.567'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +2.20%  '
$ws.Cells.Item(16, 5).Value = '  +0.18%  '
$ws.Cells.Item(17, 4).Value = '28.312.32'
$ws.Cells.Item(17, 5).Value = '  +3.96%  '
$ws.Cells.Item(18, 4).Value = '''246.28'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +3.77%  '
$ws.Cells.Item(19, 4).Value = '0.0₃0752'
$ws.Cells.Item(19, 5).Value = '  +1.01%  '
$ws.Cells.Item(20, 4).Value = '''7.92'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -2.54%  '
$ws.Cells.Item(21, 5).Value = '  -0.14%  '
$ws.Cells.Item(22, 4).Value = '''4.63'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +1.53%  '
$ws.Cells.Item(23, 4).Value = '''9.67'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.33%  '
$ws.Cells.Item(24, 4).Value = '''2.07'
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -1.91%  '
$ws.Cells.Item(25, 4).Value = '''149.15'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +1.24%  '
$ws.Cells.Item(26, 4).Value = '''7.50'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +2.46%  '
$ws.Cells.Item(27, 4).Value = '''16.67'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +1.26%  '
$ws.Cells.Item(28, 5).Value = '  +0.20%  '
$ws.Cells.Item(29, 5).Value = '  -0.03%  '
$ws.Cells.Item(30, 5).Value = '  +2.82%  '
$ws.Cells.Item(31, 5).Value = '  +2.66%  '
$ws.Cells.Item(32, 4).Value = '''3.43'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +0.74%  '
$ws.Cells.Item(33, 2).Value = 'Maker'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(33, 4).Value = '1.486.57'
$ws.Cells.Item(33, 5).Value = '  -4.24%  '
$ws.Cells.Item(34, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(34, 4).Value = '''3.26'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +0.64%  '
$ws.Cells.Item(35, 5).Value = '  -2.44%  '
$ws.Cells.Item(36, 4).Value = '''0.983'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +3.61%  '
$ws.Cells.Item(37, 5).Value = '  +0.05%  '
$ws.Cells.Item(38, 5).Value = '  -0.09%  '
$ws.Cells.Item(39, 5).Value = '  +1.17%  '
$ws.Cells.Item(40, 5).Value = '  +0.31%  '
$ws.Cells.Item(41, 4).Value = '''69.94'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +0.83%  '
$ws.Cells.Item(42, 5).Value = '  -0.24%  '
$ws.Cells.Item(44, 4).Value = '1.879.83'
$ws.Cells.Item(44, 5).Value = '  +2.38%  '
$ws.Cells.Item(45, 5).Value = '  +1.25%  '
$ws.Cells.Item(46, 4).Value = '''0.803'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +1.62%  '
$ws.Cells.Item(47, 4).Value = '''1.74'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +7.82%  '
$ws.Cells.Item(48, 2).Value = 'Quant'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(48, 4).Value = '''90.38'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.88%  '
$ws.Cells.Item(49, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(49, 4).Value = '0.0₆0113'
$ws.Cells.Item(49, 5).Value = '  +3.07%  '
$ws.Cells.Item(50, 4).Value = '''8.19'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -1.18%  '
$ws.Cells.Item(51, 5).Value = '  -0.49%  '
